$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 696
$ws.Range("A29").Value = 733
$ws.Range("A30").Value = 750
$ws.Range("K30").Value = 2.49
$ws.Range("A31").Value = 804
$ws.Range("A32").Value = 840
$ws.Range("A33").Value = 849
$ws.Range("A34").Value = 859
$ws.Range("A35").Value = 871
$ws.Range("A36").Value = 904
$ws.Range("A37").Value = 918
$ws.Range("A38").Value = 971
$ws.Range("A39").Value = 988
$ws.Range("A40").Value = 990
$ws.Range("A41").Value = 1009
$ws.Range("A42").Value = 1114
$ws.Range("A43").Value = 1196
$ws.Range("A44").Value = 1200
$ws.Range("A45").Value = 1201
$ws.Range("A46").Value = 1241
$ws.Range("A47").Value = 1245
$ws.Range("A48").Value = 1310
$ws.Range("A49").Value = 1312
$ws.Range("A50").Value = 1388
$ws.Range("A51").Value = 1404
$ws.Range("K51").Value = 3.79
$ws.Range("A52").Value = 1424
$ws.Range("A53").Value = 1429
$ws.Range("A54").Value = 1434
$ws.Range("K54").Value = 11.02
$ws.Range("A55").Value = 1437
$ws.Range("K55").Value = 10.88
$ws.Range("A56").Value = 1443
$ws.Range("A57").Value = 1560
$ws.Range("A58").Value = 1578
$ws.Range("A59").Value = 1602
$ws.Range("K59").Value = 6.42
$ws.Range("A60").Value = 1647
$ws.Range("A61").Value = 1666
$ws.Range("A62").Value = 1770
$ws.Range("A63").Value = 1819
$ws.Range("A64").Value = 1849
$ws.Range("A65").Value = 1895
$ws.Range("A66").Value = 1910
$ws.Range("A67").Value = 1913
$ws.Range("K67").Value = 9.23
$ws.Range("A68").Value = 1929
$ws.Range("A69").Value = 1983
$ws.Range("A70").Value = 2029
$ws.Range("A71").Value = 2034
$ws.Range("A72").Value = 2039
$ws.Range("A73").Value = 2042
$ws.Range("A74").Value = 2078
$ws.Range("K74").Value = 9.51
$ws.Range("A75").Value = 2207
$ws.Range("A76").Value = 2266
$ws.Range("A77").Value = 2274
$ws.Range("A78").Value = 2334
$ws.Range("A79").Value = 2377
$ws.Range("A80").Value = 2423
$ws.Range("A81").Value = 2551
$ws.Range("A82").Value = 2564
$ws.Range("A83").Value = 2622
$ws.Range("K83").Value = 5.6
$ws.Range("A84").Value = 2681
$ws.Range("A85").Value = 2684
$ws.Range("K85").Value = 2.64
$ws.Range("A86").Value = 2755
$ws.Range("A87").Value = 2800
$ws.Range("A88").Value = 2937
$ws.Range("A89").Value = 2968
$ws.Range("A90").Value = 2979
$ws.Range("A91").Value = 2982
$ws.Range("A92").Value = 3097
$ws.Range("A93").Value = 3112
$ws.Range("A94").Value = 3154
$ws.Range("A95").Value = 3215
$ws.Range("A96").Value = 3286
$ws.Range("A97").Value = 3356
$ws.Range("A98").Value = 3361
$ws.Range("A99").Value = 3409
$ws.Range("A100").Value = 3581
$ws.Range("A101").Value = 3589
$ws.Range("A102").Value = 3605
$ws.Range("A103").Value = 3638
$ws.Range("A104").Value = 3649
$ws.Range("A105").Value = 3764
$ws.Range("A106").Value = 3791
$ws.Range("A107").Value = 3830
$ws.Range("A108").Value = 3846
$ws.Range("A109").Value = 3940
$ws.Range("A110").Value = 3987
$ws.Range("A111").Value = 3995
$ws.Range("A112").Value = 3996
$ws.Range("A113").Value = 4017
$ws.Range("A114").Value = 4029
$ws.Range("A115").Value = 4050
$ws.Range("A116").Value = 4140
$ws.Range("A117").Value = 4174
$ws.Range("A118").Value = 4227
$ws.Range("A119").Value = 4266
$ws.Range("A120").Value = 4271
$ws.Range("A121").Value = 4317
$ws.Range("A122").Value = 4351
$ws.Range("A123").Value = 4483
$ws.Range("A124").Value = 4509
$ws.Range("A125").Value = 4537
$ws.Range("A126").Value = 4582
$ws.Range("A127").Value = 4612
$ws.Range("A128").Value = 4627
$ws.Range("K128").Value = 3.41
$ws.Range("A129").Value = 4646
$ws.Range("A130").Value = 4710
$ws.Range("A131").Value = 4823
$ws.Range("A132").Value = 4825
$ws.Range("A133").Value = 4839
$ws.Range("A134").Value = 4865
$ws.Range("A135").Value = 4870
$ws.Range("A136").Value = 4908
$ws.Range("A137").Value = 4931
$ws.Range("A138").Value = 4948
$ws.Range("A139").Value = 5058
$ws.Range("K139").Value = 8.48
$ws.Range("A140").Value = 5073
$ws.Range("A141").Value = 5087
$ws.Range("A142").Value = 5109
$ws.Range("A143").Value = 5112
$ws.Range("A144").Value = 5125
$ws.Range("A145").Value = 5169
$ws.Range("A146").Value = 5226
$ws.Range("A147").Value = 5231
$ws.Range("A148").Value = 5253
$ws.Range("A149").Value = 5283
